$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.250.64"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.519.08"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.38"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.93"
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.516"
$ws.Range("E8").Value = "  +1.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.519.06"
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").Value = "  +1.61%  "
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("E12").Value = "  +4.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.90"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.982.61"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.123.82"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.83"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.520.50"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.34"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.64"
$ws.Range("E20").Value = "  +1.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.12"
$ws.Range("E21").Value = "  -1.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.91"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.00"
$ws.Range("E23").Value = "  +1.71%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  +2.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.96"
$ws.Range("E26").Value = "  -2.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.89"
$ws.Range("E27").Value = "  -3.39%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0894"
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.82"
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "462.49"
$ws.Range("E32").Value = "  -3.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.23"
$ws.Range("E33").Value = "  -3.83%  "
$ws.Range("E34").Value = "  -1.03%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.116"
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.49"
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.55"
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.319"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.71"
$ws.Range("E42").Value = "  +0.56%  "
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.14"
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.13"
$ws.Range("E45").Value = "  -13.02%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.25"
$ws.Range("E46").Value = "  -4.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.78"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.522"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.47"
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0730"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.580"
$ws.Range("E51").Value = "  -2.64%  "
